$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptocurrency price/volume table (columns B-E, rows 2-51)
# with the latest scraped values. Every updated cell is forced to the
# "Text" number format before the value is written so that strings which
# look numeric (e.g. "231.65", "0.0923") are kept as literal text instead
# of being parsed into floating point numbers by Excel. The style is then
# reset back to "Normal" so no extra formatting is left behind on the cell.
function Set-TextValue($address, $value) {
    $rng = $ws.Range($address)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "37.106.36"
Set-TextValue "E2" "  +1.50%  "
Set-TextValue "D3" "2.048.13"
Set-TextValue "E3" "  +0.85%  "
Set-TextValue "E4" "  -0.04%  "
Set-TextValue "D5" "231.65"
Set-TextValue "E5" "  -0.34%  "
Set-TextValue "D6" "0.617"
Set-TextValue "E6" "  +3.04%  "
Set-TextValue "E7" "  -0.02%  "
Set-TextValue "D8" "57.21"
Set-TextValue "E8" "  +3.72%  "
Set-TextValue "E9" "  +2.58%  "
Set-TextValue "D10" "57.59"
Set-TextValue "E10" "  +0.77%  "
Set-TextValue "E11" "  +1.03%  "
Set-TextValue "E12" "  +0.90%  "
Set-TextValue "D13" "2.352.90"
Set-TextValue "E13" "  +0.95%  "
Set-TextValue "D14" "14.24"
Set-TextValue "E14" "  -0.64%  "
Set-TextValue "D15" "20.76"
Set-TextValue "E15" "  +3.06%  "
Set-TextValue "E16" "  +0.98%  "
Set-TextValue "E17" "  +0.19%  "
Set-TextValue "D18" "2.046.31"
Set-TextValue "E18" "  +1.70%  "
Set-TextValue "D19" "37.033.81"
Set-TextValue "E19" "  +0.75%  "
Set-TextValue "D20" "6.29"
Set-TextValue "E20" "  +13.16%  "
Set-TextValue "D21" "68.75"
Set-TextValue "E21" "  +1.80%  "
Set-TextValue "E22" "  +1.26%  "
Set-TextValue "D23" "224.19"
Set-TextValue "E23" "  +1.49%  "
Set-TextValue "E24" "  +0.02%  "
Set-TextValue "D25" "2.42"
Set-TextValue "E25" "  +1.45%  "
Set-TextValue "E26" "  +0.02%  "
Set-TextValue "D27" "165.30"
Set-TextValue "E27" "  +1.48%  "
Set-TextValue "E28" "  +7.42%  "
Set-TextValue "E29" "  +0.84%  "
Set-TextValue "D30" "18.98"
Set-TextValue "E30" "  +0.55%  "
Set-TextValue "D31" "0.124"
Set-TextValue "E31" "  -3.04%  "
Set-TextValue "D32" "0.116"
Set-TextValue "E32" "  -0.61%  "
Set-TextValue "D33" "4.45"
Set-TextValue "E33" "  +1.90%  "
Set-TextValue "E34" "  +2.16%  "
Set-TextValue "E35" "  +1.24%  "
Set-TextValue "D36" "4.51"
Set-TextValue "E36" "  +5.63%  "
Set-TextValue "E37" "  -0.03%  "
Set-TextValue "E38" "  -0.72%  "
Set-TextValue "E39" "  -1.21%  "
Set-TextValue "E40" "  -2.33%  "
Set-TextValue "D41" "4.51"
Set-TextValue "E41" "  +9.04%  "
Set-TextValue "E42" "  +1.29%  "
Set-TextValue "D43" "1.481.52"
Set-TextValue "E43" "  +0.54%  "
Set-TextValue "D44" "95.90"
Set-TextValue "E44" "  +3.14%  "
Set-TextValue "B45" "TrustWalletToken"
Set-TextValue "C45" "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
Set-TextValue "D45" "1.16"
Set-TextValue "E45" "  +3.39%  "
Set-TextValue "B46" "Cronos"
Set-TextValue "C46" "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue "D46" "0.0923"
Set-TextValue "E46" "  -0.45%  "
Set-TextValue "E47" "  +2.87%  "
Set-TextValue "D48" "15.19"
Set-TextValue "E48" "  -3.18%  "
Set-TextValue "D49" "1.01"
Set-TextValue "E49" "  +0.98%  "
Set-TextValue "E50" "  +2.98%  "
Set-TextValue "E51" "  +1.02%  "
